$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "phi3:14b-medium-4k-instruct-q5_K_M"
$ws.Range("B12").Value = "llama3:70b"
$ws.Range("C12").Value = 42
$ws.Range("D12").Value = 200
$ws.Range("E12").Value = 2846.55
$ws.Range("F12").Value = 306.1151
$ws.Range("G12").Value = 2.5
$ws.Range("H12").Value = "phi3_14b_medium_4k_instruct_q5_K_M_llama3_70b_42_200_val.txt"
$ws.Range("I12").Value = 477.93
$ws.Range("J12").Value = 18.75
$ws.Range("K12").Value = "phi3_14b_medium_4k_instruct_q5_K_M_llama3_70b_42_200_test.txt"
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 4
$ws.Range("N12").Value = 148.71
$ws.Range("O12").Value = 57.5
$ws.Range("P12").Value = "phi3_14b_medium_4k_instruct_q5_K_M_llama3_70b_42_200_val_labeled.txt"
$ws.Range("Q12").Value = 276.76
$ws.Range("R12").Value = 57.5
$ws.Range("S12").Value = "phi3_14b_medium_4k_instruct_q5_K_M_llama3_70b_42_200_test_labeled.txt"
$ws.Range("T12").Value = 1168.15
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 2
$ws.Range("W12").Value = 133.24
$ws.Range("X12").Value = 72.5
$ws.Range("Y12").Value = "phi3_14b_medium_4k_instruct_q5_K_M_llama3_70b_42_200_val_bootstrap.txt"
$ws.Range("Z12").Value = 335.65
$ws.Range("AA12").Value = 52.5
$ws.Range("AB12").Value = "phi3_14b_medium_4k_instruct_q5_K_M_llama3_70b_42_200_test_bootstrap.txt"
